# Update crypto price/volume data per GitHub Actions scrape run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Keep these cells stored as text (matching the original inline-string cells)
# instead of letting Excel auto-convert numeric-looking text into numbers/percentages.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "328.76"
$ws.Range("E2").Value = "1.68%"
$ws.Range("D3").Value = "41.31"
$ws.Range("E3").Value = "4.90%"
$ws.Range("D4").Value = "5.641"
$ws.Range("E4").Value = "-1.09%"
$ws.Range("D5").Value = "0.08172"
$ws.Range("E5").Value = "2.23%"
$ws.Range("E6").Value = "1.63%"
$ws.Range("D7").Value = "2.012"
$ws.Range("E7").Value = "1.53%"
$ws.Range("D8").Value = "4.496"
$ws.Range("E8").Value = "-0.80%"
$ws.Range("E9").Value = "1.20%"
$ws.Range("D10").Value = "0.9200"
$ws.Range("E10").Value = "-0.70%"
$ws.Range("D11").Value = "0.1273"
$ws.Range("E11").Value = "2.14%"
$ws.Range("D12").Value = "0.1966"
$ws.Range("E12").Value = "0.45%"
$ws.Range("D13").Value = "0.09440"
$ws.Range("E13").Value = "2.68%"
$ws.Range("D14").Value = "0.03753"
$ws.Range("E14").Value = "2.90%"
$ws.Range("D15").Value = "0.1060"
$ws.Range("E15").Value = "1.02%"
$ws.Range("E16").Value = "2.09%"
$ws.Range("D17").Value = "0.006133"
$ws.Range("E17").Value = "-1.93%"
$ws.Range("D19").Value = "3.446"
$ws.Range("E19").Value = "2.81%"
$ws.Range("E20").Value = "-1.42%"
$ws.Range("D21").Value = "8.363"
$ws.Range("E21").Value = "-3.86%"
$ws.Range("E22").Value = "0.02%"
$ws.Range("D23").Value = "0.2412"
$ws.Range("E23").Value = "-1.55%"
$ws.Range("D24").Value = "0.04394"
$ws.Range("E24").Value = "-0.49%"
$ws.Range("D25").Value = "0.001257"
$ws.Range("E25").Value = "-0.52%"
$ws.Range("D26").Value = "0.004320"
$ws.Range("E26").Value = "-4.39%"
$ws.Range("E27").Value = "4.35%"
$ws.Range("D39").Value = "0.02803"
$ws.Range("E39").Value = "11.32%"
$ws.Range("D40").Value = "0.05409"
$ws.Range("D41").Value = "0.007681"
$ws.Range("E41").Value = "3.58%"
$ws.Range("D42").Value = "0.1416"
$ws.Range("E42").Value = "0.95%"
$ws.Range("D43").Value = "0.008947"
$ws.Range("E43").Value = "-6.34%"
$ws.Range("E44").Value = "2.54%"
$ws.Range("D45").Value = "0.01154"
$ws.Range("E45").Value = "7.86%"
$ws.Range("D46").Value = "0.00006601"
$ws.Range("E46").Value = "-2.57%"
$ws.Range("E47").Value = "-0.02%"
$ws.Range("E48").Value = "7.57%"
$ws.Range("D49").Value = "0.002282"
$ws.Range("E49").Value = "-0.50%"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "-0.02%"
